$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'76.072.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.96%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.862.71"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +7.39%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.08%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'195.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +4.38%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'597.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.77%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.07%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.552"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +3.46%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.192"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.40%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'2.864.53"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +7.51%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.390"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +9.18%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -1.99%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'4.90"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +3.58%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'3.394.30"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +7.34%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'76.175.70"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.64%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'27.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +3.53%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.0000188"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.71%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.878.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +7.62%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'9.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.80%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'12.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +4.80%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'381.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +2.59%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'2.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +2.53%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'4.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.90%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'71.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +2.54%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +0.08%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'3.030.06"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +8.41%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'4.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.37%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'9.69"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +3.79%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.0000104"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +9.85%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -0.19%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'1.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -1.92%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'507.80"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -2.36%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'7.70"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.40%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.80"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +2.75%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.04%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'167.19"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +2.26%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'19.93"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +3.93%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.117"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.55%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'19.52"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.82%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'184.63"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +8.60%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -0.06%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.343"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +4.05%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'5.03"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.96%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.94%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.0916"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +8.30%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'1.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +2.82%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'40.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +2.74%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.28%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.575"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +7.65%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.675"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +14.29%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'3.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +2.64%  "
$ws.Range("E51").Style = "Normal"
